# Rename the header/footer logo pictures so their InlineShape.Name
# (the underlying <wp:docPr name="..."/> attribute) changes:
#   - the two Pearson logo pictures (footers):  image2.png -> image1.png
#   - the BTEC logo picture (header):           image1.jpg -> image2.jpg
#
# The pictures are matched by their (stable) AlternativeText rather than a
# blind positional index, so the script is resilient to header/footer
# ordering.

$d = $word.ActiveDocument

$pearsonAlt = "Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png"
$btecAlt = "BTec_Logo-Orange"

for ($s = 1; $s -le $d.Sections.Count; $s++) {
    $section = $d.Sections.Item($s)

    for ($h = 1; $h -le $section.Headers.Count; $h++) {
        $header = $section.Headers.Item($h)
        if ($header.Exists) {
            $shapes = $header.Range.InlineShapes
            for ($i = 1; $i -le $shapes.Count; $i++) {
                $shape = $shapes.Item($i)
                if ($shape.AlternativeText -eq $btecAlt) {
                    $shape.Name = "image2.jpg"
                }
            }
        }
    }

    for ($f = 1; $f -le $section.Footers.Count; $f++) {
        $footer = $section.Footers.Item($f)
        if ($footer.Exists) {
            $shapes = $footer.Range.InlineShapes
            for ($i = 1; $i -le $shapes.Count; $i++) {
                $shape = $shapes.Item($i)
                if ($shape.AlternativeText -eq $pearsonAlt) {
                    $shape.Name = "image1.png"
                }
            }
        }
    }
}
